# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
# (Updated cryptos list on Mon Apr 29 12:56:11 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.670.61'
$ws.Range('E2').Value = '  -1.48%  '
# Row 3
$ws.Range('D3').Value = '3.170.29'
$ws.Range('E3').Value = '  -3.68%  '
# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.02%  '
# Row 5
$ws.Range('D5').Value = '''589.41'
$ws.Range('E5').Value = '  -1.86%  '
# Row 6
$ws.Range('D6').Value = '''135.54'
$ws.Range('E6').Value = '  -4.36%  '
# Row 7
$ws.Range('E7').Value = '  +0.02%  '
# Row 8
$ws.Range('D8').Value = '3.165.64'
$ws.Range('E8').Value = '  -3.79%  '
# Row 9
$ws.Range('D9').Value = '''0.508'
$ws.Range('E9').Value = '  -2.11%  '
# Row 10
$ws.Range('E10').Value = '  -5.13%  '
# Row 11
$ws.Range('E11').Value = '  -3.04%  '
# Row 12
$ws.Range('E12').Value = '  -3.07%  '
# Row 13
$ws.Range('E13').Value = '  -5.06%  '
# Row 14
$ws.Range('D14').Value = '''33.45'
$ws.Range('E14').Value = '  -3.00%  '
# Row 15
$ws.Range('D15').Value = '3.691.09'
$ws.Range('E15').Value = '  -3.78%  '
# Row 16
$ws.Range('E16').Value = '  -2.28%  '
# Row 17
$ws.Range('D17').Value = '3.166.72'
$ws.Range('E17').Value = '  -3.73%  '
# Row 18
$ws.Range('D18').Value = '62.617.90'
$ws.Range('E18').Value = '  -1.62%  '
# Row 19
$ws.Range('D19').Value = '''6.54'
$ws.Range('E19').Value = '  -4.17%  '
# Row 20
$ws.Range('D20').Value = '''454.12'
$ws.Range('E20').Value = '  -4.98%  '
# Row 21
$ws.Range('D21').Value = '''13.96'
$ws.Range('E21').Value = '  -0.85%  '
# Row 22
$ws.Range('E22').Value = '  -3.61%  '
# Row 23
$ws.Range('D23').Value = '''7.61'
$ws.Range('E23').Value = '  -5.39%  '
# Row 24
$ws.Range('D24').Value = '''13.39'
$ws.Range('E24').Value = '  -1.63%  '
# Row 25
$ws.Range('D25').Value = '''83.65'
$ws.Range('E25').Value = '  -0.62%  '
# Row 26
$ws.Range('E26').Value = '  -0.02%  '
# Row 27
$ws.Range('D27').Value = '''2.69'
$ws.Range('E27').Value = '  -2.40%  '
# Row 28
$ws.Range('D28').Value = '''0.999'
$ws.Range('E28').Value = '  -0.03%  '
# Row 29
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '''7.72'
$ws.Range('E29').Value = '  -4.39%  '
# Row 30
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '''6.77'
$ws.Range('E30').Value = '  -6.19%  '
# Row 31
$ws.Range('E31').Value = '  -5.97%  '
# Row 32
$ws.Range('D32').Value = '''27.28'
$ws.Range('E32').Value = '  -4.65%  '
# Row 33
$ws.Range('E33').Value = '  -1.11%  '
# Row 34
$ws.Range('D34').Value = '''2.39'
$ws.Range('E34').Value = '  -5.79%  '
# Row 35
$ws.Range('E35').Value = '  -6.66%  '
# Row 36
$ws.Range('D36').Value = '''5.87'
$ws.Range('E36').Value = '  -1.71%  '
# Row 37
$ws.Range('D37').Value = '''51.19'
$ws.Range('E37').Value = '  -3.72%  '
# Row 38
$ws.Range('D38').Value = '0.0₃0702'
$ws.Range('E38').Value = '  -4.56%  '
# Row 39
$ws.Range('D39').Value = '''0.0386'
$ws.Range('E39').Value = '  -2.85%  '
# Row 40
$ws.Range('D40').Value = '''2.70'
$ws.Range('E40').Value = '  -0.98%  '
# Row 41
$ws.Range('D41').Value = '''402.50'
$ws.Range('E41').Value = '  -5.81%  '
# Row 42
$ws.Range('D42').Value = '''8.01'
$ws.Range('E42').Value = '  -3.78%  '
# Row 43
$ws.Range('E43').Value = '  -0.98%  '
# Row 44
$ws.Range('D44').Value = '2.798.68'
$ws.Range('E44').Value = '  -8.87%  '
# Row 45
$ws.Range('D45').Value = '''0.250'
$ws.Range('E45').Value = '  -5.34%  '
# Row 46
$ws.Range('D46').Value = '''2.14'
$ws.Range('E46').Value = '  -2.20%  '
# Row 48
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''125.58'
$ws.Range('E48').Value = '  -1.29%  '
# Row 49
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = '''34.89'
$ws.Range('E49').Value = '  -1.63%  '
# Row 50
$ws.Range('D50').Value = '''25.35'
$ws.Range('E50').Value = '  -3.28%  '
# Row 51
$ws.Range('E51').Value = '  -3.13%  '
